# Data Mining Project.pptx - "Added more to presentation outline"
#
# Before: 6 slides
#   1. Data Mining Project (title slide)                         [SlideID 256]
#   2. (no title) Memory-based Collaborative filtering / KNN     [SlideID 257]
#   3. (no title) Model Based Collaborative filtering / SVD ...  [SlideID 258]
#   4. (no title, empty content)                                 [SlideID 259]
#   5. (no title, empty content)                                 [SlideID 260]
#   6. Sources                                                   [SlideID 261]
#
# After: 11 slides (5 new slides inserted, 3 existing slides retitled/filled)
#   1. Data Mining Project (title slide)               [256] (unchanged)
#   2. Outline                                         [262] (NEW)
#   3. Introduction to project                         [263] (NEW)
#   4. Purpose                                         [264] (NEW)
#   5. Dataset overview                                [265] (NEW)
#   6. Objective 1 - predict using ....                [257] (title added, content kept)
#   7. Objective 1 - Conclusion                        [266] (NEW)
#   8. Objective 2 - predict using ....                [258] (title added, content kept)
#   9. Objective 2 - Conclusion                        [259] (title + content added)
#  10. (still blank)                                   [260] (unchanged)
#  11. Sources                                         [261] (unchanged)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Insert the five brand-new slides at the right slots (layout 2 = "Title
#    and Content", the same layout every other content slide in this deck
#    uses). Inserting in ascending index order reproduces the exact SlideID
#    allocation order (262, 263, 264, 265, 266) and final slide order.
# ---------------------------------------------------------------------------

$sOutline = $p.Slides.Add(2, 2)
$sOutline.Shapes.Item(1).TextFrame.TextRange.Text = "Outline "

$sIntro = $p.Slides.Add(3, 2)
$sIntro.Shapes.Item(1).TextFrame.TextRange.Text = "Introduction to project"
$sIntro.Shapes.Item(2).TextFrame.TextRange.Text = "I’m going to ….."

$sPurpose = $p.Slides.Add(4, 2)
$sPurpose.Shapes.Item(1).TextFrame.TextRange.Text = "Purpose"
$sPurpose.Shapes.Item(2).TextFrame.TextRange.Text = "Statement of purpose`r“To predict user movie ratings using this and that algorithm”`rPredict movies???"

$sDataset = $p.Slides.Add(5, 2)
$sDataset.Shapes.Item(1).TextFrame.TextRange.Text = "Dataset overview"
$sDataset.Shapes.Item(2).TextFrame.TextRange.Text = "Acquired`rInstances, attributes`rNo need to clean MovieLens"

$sObj1Concl = $p.Slides.Add(7, 2)
$sObj1Concl.Shapes.Item(1).TextFrame.TextRange.Text = "Objective 1 – Conclusion"
$sObj1Concl.Shapes.Item(2).TextFrame.TextRange.Text = "evaluation"

# ---------------------------------------------------------------------------
# 2. Retitle the pre-existing slides that now sit at positions 6, 8 and 9.
#    Their bodies already hold the right content (untouched), except for
#    slide 9 ("Objective 2 - Conclusion"), which was completely empty and
#    now needs both a title and body text.
# ---------------------------------------------------------------------------

$sObj1 = $p.Slides.Item(6)
$sObj1.Shapes.Item(1).TextFrame.TextRange.Text = "Objective 1 – predict using …."

$sObj2 = $p.Slides.Item(8)
$sObj2.Shapes.Item(1).TextFrame.TextRange.Text = "Objective 2 – predict using …."

$sObj2Concl = $p.Slides.Item(9)
$sObj2Concl.Shapes.Item(1).TextFrame.TextRange.Text = "Objective 2 – Conclusion"
$sObj2Concl.Shapes.Item(2).TextFrame.TextRange.Text = "Evaluation"
